# chore: update Sheets via scheduled runner
# Refresh cached market-board price / profit figures (columns H-N) for the
# Leve tables on each crafting-job sheet. Only the computed price/profit
# columns change; leve metadata (A-G) is untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1416.7858
$ws.Range("I62").Value = 1319.0834
$ws.Range("K62").Value = 1319.0834
$ws.Range("M62").Value = -695.0834

$ws.Range("H65").Value = 1416.7858
$ws.Range("I65").Value = 1319.0834
$ws.Range("K65").Value = 6595.416999999999
$ws.Range("M65").Value = -3475.416999999999

$ws.Range("H132").Value = 1545.7073
$ws.Range("I132").Value = 1569.35
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 4708.049999999999
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -2178.049999999999
$ws.Range("N132").Value = -6860

$ws.Range("H135").Value = 150002260
$ws.Range("I135").Value = 71431530
$ws.Range("K135").Value = 642883770
$ws.Range("M135").Value = -642881235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5956.891
$ws.Range("I61").Value = 4038.138
$ws.Range("J61").Value = 9230.058999999999
$ws.Range("K61").Value = 4038.138
$ws.Range("L61").Value = 9230.058999999999
$ws.Range("M61").Value = -3826.138
$ws.Range("N61").Value = -9654.058999999999

$ws.Range("H74").Value = 6323.558
$ws.Range("I74").Value = 4661
$ws.Range("K74").Value = 4661
$ws.Range("M74").Value = -3787

$ws.Range("H77").Value = 6323.558
$ws.Range("I77").Value = 4661
$ws.Range("K77").Value = 23305
$ws.Range("M77").Value = -18937

$ws.Range("H102").Value = 2250.7827
$ws.Range("I102").Value = 1680.5294
$ws.Range("K102").Value = 1680.5294
$ws.Range("M102").Value = -58.5293999999999

$ws.Range("H136").Value = 5956.891
$ws.Range("I136").Value = 4038.138
$ws.Range("J136").Value = 9230.058999999999
$ws.Range("K136").Value = 12114.414
$ws.Range("L136").Value = 27690.177
$ws.Range("M136").Value = -9564.414000000001
$ws.Range("N136").Value = -32790.177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 64550
$ws.Range("J13").Value = 64550
$ws.Range("L13").Value = 64550
$ws.Range("N13").Value = -64886

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4261.222
$ws.Range("I31").Value = 4859.4814
$ws.Range("J31").Value = 3363.8333
$ws.Range("K31").Value = 4859.4814
$ws.Range("L31").Value = 3363.8333
$ws.Range("M31").Value = -4564.4814
$ws.Range("N31").Value = -3953.8333

$ws.Range("H34").Value = 4261.222
$ws.Range("I34").Value = 4859.4814
$ws.Range("J34").Value = 3363.8333
$ws.Range("K34").Value = 4859.4814
$ws.Range("L34").Value = 3363.8333
$ws.Range("M34").Value = -4657.4814
$ws.Range("N34").Value = -3767.8333

$ws.Range("H58").Value = 1717191.8
$ws.Range("I58").Value = 2332159.5
$ws.Range("J58").Value = 4067.4285
$ws.Range("K58").Value = 2332159.5
$ws.Range("L58").Value = 4067.4285
$ws.Range("M58").Value = -2331956.5
$ws.Range("N58").Value = -4473.4285

$ws.Range("H99").Value = 1982.6522
$ws.Range("I99").Value = 1622.125
$ws.Range("J99").Value = 2806.7144
$ws.Range("K99").Value = 1622.125
$ws.Range("L99").Value = 2806.7144
$ws.Range("M99").Value = -124.125
$ws.Range("N99").Value = -5802.7144

$ws.Range("H115").Value = 37040
$ws.Range("J115").Value = 37040
$ws.Range("L115").Value = 37040
$ws.Range("N115").Value = -39390

$ws.Range("H122").Value = 11504.714
$ws.Range("I122").Value = 4836.5
$ws.Range("J122").Value = 51514
$ws.Range("K122").Value = 14509.5
$ws.Range("L122").Value = 154542
$ws.Range("M122").Value = -12059.5
$ws.Range("N122").Value = -159442

$ws.Range("H126").Value = 1982.6522
$ws.Range("I126").Value = 1622.125
$ws.Range("J126").Value = 2806.7144
$ws.Range("K126").Value = 4866.375
$ws.Range("L126").Value = 8420.143199999999
$ws.Range("M126").Value = -2396.375
$ws.Range("N126").Value = -13360.1432

$ws.Range("H136").Value = 1717191.8
$ws.Range("I136").Value = 2332159.5
$ws.Range("J136").Value = 4067.4285
$ws.Range("K136").Value = 6996478.5
$ws.Range("L136").Value = 12202.2855
$ws.Range("M136").Value = -6993928.5
$ws.Range("N136").Value = -17302.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2395.4546
$ws.Range("I81").Value = 713.75
$ws.Range("J81").Value = 3356.4285
$ws.Range("K81").Value = 2141.25
$ws.Range("L81").Value = 10069.2855
$ws.Range("M81").Value = -1018.25
$ws.Range("N81").Value = -12315.2855

$ws.Range("H84").Value = 2395.4546
$ws.Range("I84").Value = 713.75
$ws.Range("J84").Value = 3356.4285
$ws.Range("K84").Value = 6423.75
$ws.Range("L84").Value = 30207.8565
$ws.Range("M84").Value = -807.75
$ws.Range("N84").Value = -41439.8565

$ws.Range("H113").Value = 684.1818
$ws.Range("I113").Value = 701.8421
$ws.Range("J113").Value = 660.2143
$ws.Range("K113").Value = 2105.5263
$ws.Range("L113").Value = 1980.6429
$ws.Range("M113").Value = 64.47370000000001
$ws.Range("N113").Value = -6320.6429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4405.6665
$ws.Range("I122").Value = 4763.077
$ws.Range("J122").Value = 3476.4
$ws.Range("K122").Value = 14289.231
$ws.Range("L122").Value = 10429.2
$ws.Range("M122").Value = -11839.231
$ws.Range("N122").Value = -15329.2

$ws.Range("H126").Value = 2976.1177
$ws.Range("I126").Value = 1999.25
$ws.Range("J126").Value = 3844.4443
$ws.Range("K126").Value = 5997.75
$ws.Range("L126").Value = 11533.3329
$ws.Range("M126").Value = -3527.75
$ws.Range("N126").Value = -16473.3329

$ws.Range("H132").Value = 2433.111
$ws.Range("I132").Value = 2226.2104
$ws.Range("J132").Value = 2924.5
$ws.Range("K132").Value = 6678.6312
$ws.Range("L132").Value = 8773.5
$ws.Range("M132").Value = -4148.6312
$ws.Range("N132").Value = -13833.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2917.5557
$ws.Range("I40").Value = 2917.5557
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2917.5557
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2781.5557
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 5655.237
$ws.Range("I132").Value = 6530.24
$ws.Range("J132").Value = 3972.5386
$ws.Range("K132").Value = 19590.72
$ws.Range("L132").Value = 11917.6158
$ws.Range("M132").Value = -17060.72
$ws.Range("N132").Value = -16977.6158

$ws.Range("H136").Value = 4609.102
$ws.Range("I136").Value = 2299.111
$ws.Range("J136").Value = 7444.091
$ws.Range("K136").Value = 6897.333
$ws.Range("L136").Value = 22332.273
$ws.Range("M136").Value = -4347.333
$ws.Range("N136").Value = -27432.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5477.6816
$ws.Range("I113").Value = 10501.6
$ws.Range("J113").Value = 1291.0834
$ws.Range("K113").Value = 31504.8
$ws.Range("L113").Value = 3873.2502
$ws.Range("M113").Value = -29334.8
$ws.Range("N113").Value = -8213.2502

$ws.Range("H118").Value = 67300
$ws.Range("J118").Value = 67300
$ws.Range("L118").Value = 67300
$ws.Range("N118").Value = -70614

$ws.Range("H126").Value = 1478.579
$ws.Range("I126").Value = 1493
$ws.Range("J126").Value = 1401.6666
$ws.Range("K126").Value = 4479
$ws.Range("L126").Value = 4204.9998
$ws.Range("M126").Value = -2009
$ws.Range("N126").Value = -9144.9998

$ws.Range("H132").Value = 2054.973
$ws.Range("I132").Value = 1002.875
$ws.Range("J132").Value = 2856.5715
$ws.Range("K132").Value = 3008.625
$ws.Range("L132").Value = 8569.7145
$ws.Range("M132").Value = -478.625
$ws.Range("N132").Value = -13629.7145

$ws.Range("H133").Value = 54290.5
$ws.Range("J133").Value = 54290.5
$ws.Range("L133").Value = 54290.5
$ws.Range("N133").Value = -64410.5
